$wb = $excel.ActiveWorkbook

# --- Sheet ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H43").Value = 445
$ws.Range("I43").Value = 0
$ws.Range("J43").Value = 445
$ws.Range("K43").Value = 0
$ws.Range("L43").Value = 445
$ws.Range("M43").ClearContents()
$ws.Range("N43").Value = -583
$ws.Range("H76").Value = 3089633.8
$ws.Range("I76").Value = 3971289.8
$ws.Range("J76").Value = 3837.5
$ws.Range("K76").Value = 3971289.8
$ws.Range("L76").Value = 3837.5
$ws.Range("M76").Value = -3970974.8
$ws.Range("N76").Value = -4467.5
$ws.Range("H79").Value = 3089633.8
$ws.Range("I79").Value = 3971289.8
$ws.Range("J79").Value = 3837.5
$ws.Range("K79").Value = 3971289.8
$ws.Range("L79").Value = 3837.5
$ws.Range("M79").Value = -3970197.8
$ws.Range("N79").Value = -6021.5
$ws.Range("H137").Value = 43479540
$ws.Range("I137").Value = 66667750
$ws.Range("J137").Value = 1648.25
$ws.Range("K137").Value = 200003250
$ws.Range("L137").Value = 4944.75
$ws.Range("M137").Value = -200000700
$ws.Range("N137").Value = -10044.75
$ws.Range("H138").Value = 7893735.5
$ws.Range("I138").Value = 5051704
$ws.Range("J138").Value = 8549589
$ws.Range("K138").Value = 15155112
$ws.Range("L138").Value = 25648767
$ws.Range("M138").Value = -15149972
$ws.Range("N138").Value = -25659047

# --- Sheet ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H101").Value = 195401.33
$ws.Range("J101").Value = 195401.33
$ws.Range("L101").Value = 195401.33
$ws.Range("N101").Value = -201891.33
$ws.Range("J102").Value = 0
$ws.Range("L102").Value = 0
$ws.Range("N102").ClearContents()
$ws.Range("H122").Value = 1872.7273
$ws.Range("I122").Value = 1775.6471
$ws.Range("J122").Value = 2202.8
$ws.Range("K122").Value = 5326.9413
$ws.Range("L122").Value = 6608.400000000001
$ws.Range("M122").Value = -2876.9413
$ws.Range("N122").Value = -11508.4
$ws.Range("H125").Value = 30000
$ws.Range("J125").Value = 0
$ws.Range("L125").Value = 0
$ws.Range("N125").ClearContents()

# --- Sheet BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H5").Value = 17744.666
$ws.Range("I5").Value = 17744.666
$ws.Range("K5").Value = 17744.666
$ws.Range("M5").Value = -17631.666
$ws.Range("H103").Value = 23922.834
$ws.Range("J103").Value = 23922.834
$ws.Range("L103").Value = 23922.834
$ws.Range("N103").Value = -26266.834
$ws.Range("H105").Value = 297445.7
$ws.Range("I105").Value = 3257.7368
$ws.Range("J105").Value = 670083.75
$ws.Range("K105").Value = 3257.7368
$ws.Range("L105").Value = 670083.75
$ws.Range("M105").Value = -1510.7368
$ws.Range("N105").Value = -673577.75

# --- Sheet CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4124.159
$ws.Range("I31").Value = 1093.0625
$ws.Range("J31").Value = 12207.083
$ws.Range("K31").Value = 1093.0625
$ws.Range("L31").Value = 12207.083
$ws.Range("M31").Value = -798.0625
$ws.Range("N31").Value = -12797.083
$ws.Range("H34").Value = 4124.159
$ws.Range("I34").Value = 1093.0625
$ws.Range("J34").Value = 12207.083
$ws.Range("K34").Value = 1093.0625
$ws.Range("L34").Value = 12207.083
$ws.Range("M34").Value = -891.0625
$ws.Range("N34").Value = -12611.083
$ws.Range("H132").Value = 2787.3784
$ws.Range("I132").Value = 2368.1614
$ws.Range("K132").Value = 7104.4842
$ws.Range("M132").Value = -4574.4842
$ws.Range("H134").Value = 2920.303
$ws.Range("I134").Value = 1533.826
$ws.Range("K134").Value = 4601.478
$ws.Range("M134").Value = -2066.478

# --- Sheet CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 14440.333
$ws.Range("I4").Value = 99.8125
$ws.Range("J4").Value = 60330
$ws.Range("K4").Value = 299.4375
$ws.Range("L4").Value = 180990
$ws.Range("M4").Value = -187.4375
$ws.Range("N4").Value = -181214
$ws.Range("H113").Value = 690.55884
$ws.Range("I113").Value = 691.4375
$ws.Range("J113").Value = 689.7778
$ws.Range("K113").Value = 2074.3125
$ws.Range("L113").Value = 2069.3334
$ws.Range("M113").Value = 95.6875
$ws.Range("N113").Value = -6409.3334
$ws.Range("H137").Value = 5945601.5
$ws.Range("J137").Value = 95278.73
$ws.Range("L137").Value = 285836.19
$ws.Range("N137").Value = -296036.19

# --- Sheet GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H4").Value = 50000
$ws.Range("J4").Value = 50000
$ws.Range("L4").Value = 50000
$ws.Range("N4").Value = -50224
$ws.Range("H70").Value = 5471.026
$ws.Range("I70").Value = 5275.6665
$ws.Range("J70").Value = 6122.222
$ws.Range("K70").Value = 5275.6665
$ws.Range("L70").Value = 6122.222
$ws.Range("M70").Value = -5005.6665
$ws.Range("N70").Value = -6662.222
$ws.Range("H73").Value = 5471.026
$ws.Range("I73").Value = 5275.6665
$ws.Range("J73").Value = 6122.222
$ws.Range("K73").Value = 5275.6665
$ws.Range("L73").Value = 6122.222
$ws.Range("M73").Value = -4339.6665
$ws.Range("N73").Value = -7994.222
$ws.Range("H80").Value = 2527
$ws.Range("I80").Value = 2355.625
$ws.Range("K80").Value = 2355.625
$ws.Range("M80").Value = -1357.625
$ws.Range("H83").Value = 2527
$ws.Range("I83").Value = 2355.625
$ws.Range("K83").Value = 11778.125
$ws.Range("M83").Value = -6786.125
$ws.Range("H97").Value = 0
$ws.Range("I97").Value = 0
$ws.Range("J97").Value = 0
$ws.Range("K97").Value = 0
$ws.Range("L97").Value = 0
$ws.Range("M97").ClearContents()
$ws.Range("N97").ClearContents()
$ws.Range("H122").Value = 465083.84
$ws.Range("I122").Value = 795334.1
$ws.Range("J122").Value = 2733.4
$ws.Range("K122").Value = 2386002.3
$ws.Range("L122").Value = 8200.200000000001
$ws.Range("M122").Value = -2383552.3
$ws.Range("N122").Value = -13100.2
$ws.Range("H123").Value = 9943.529
$ws.Range("J123").Value = 9943.529
$ws.Range("L123").Value = 9943.529
$ws.Range("N123").Value = -14843.529

# --- Sheet LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H104").Value = 33334
$ws.Range("J104").Value = 33334
$ws.Range("L104").Value = 33334
$ws.Range("N104").Value = -40322

# --- Sheet WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H123").Value = 31500
$ws.Range("J123").Value = 31500
$ws.Range("L123").Value = 31500
$ws.Range("N123").Value = -41300
$ws.Range("H132").Value = 3657.9092
$ws.Range("I132").Value = 4399.8
$ws.Range("J132").Value = 3039.6667
$ws.Range("K132").Value = 13199.4
$ws.Range("L132").Value = 9119.000100000001
$ws.Range("M132").Value = -10669.4
$ws.Range("N132").Value = -14179.0001

